$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (added in shared-string order: linklist, merge..., doing the practice)
$ws.Range("B9").Value = "linklist"
$ws.Range("D9").Value = "merge 2 linklists without using the third linklist"
$ws.Range("C1").Value = "doing the practice"

# Highlight fill on B1 (green, Accent6, Darker 25%) - no text value in B1 itself
$ws.Range("B1").Interior.ThemeColor = 10
$ws.Range("B1").Interior.TintAndShade = -0.249977111117893

# Update the active selection, matching the final cursor position recorded in the workbook
$ws.Range("J11").Select() | Out-Null
